# Update the two-digit multiplication problems in the document's tables.
# Each "old" equation text is unique in the document, so a simple
# Find/Replace (wdReplaceAll) for each pair is safe and order-independent.

$d = $word.ActiveDocument

$pairs = @(
    @("57×89=", "56×72="),
    @("50×51=", "37×46="),
    @("13×42=", "98×36="),
    @("18×84=", "90×70="),
    @("42×87=", "44×92="),
    @("51×74=", "27×12="),
    @("56×95=", "88×64="),
    @("16×87=", "23×15="),
    @("47×43=", "48×36="),
    @("14×51=", "20×62="),
    @("37×66=", "58×50="),
    @("30×23=", "53×84="),
    @("73×63=", "97×89="),
    @("41×77=", "40×36="),
    @("79×32=", "67×96="),
    @("98×55=", "11×72="),
    @("31×17=", "70×45="),
    @("97×31=", "32×46="),
    @("71×53=", "34×47="),
    @("23×76=", "46×58="),
    @("31×28=", "55×50="),
    @("40×25=", "56×81="),
    @("92×31=", "38×50="),
    @("67×20=", "59×37="),
    @("68×48=", "31×40=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
